# Update the convention/event listing data on sheets "展览" and "全部类型".
# Both sheets share an identical layout (header row 1, data rows 2-7) and
# receive identical updates: event names shift/rotate, attendee counts (F)
# are bumped slightly, and lowest-ticket-price (G) switches from a text
# percentage-like string (e.g. "54") to a numeric value scaled by 100
# (e.g. 5400).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2
    $ws.Range("C2").Value = "南宁·AP动漫游戏嘉年华"
    $ws.Range("F2").Value = 2079
    $ws.Range("G2").Value = 5400

    # Row 3
    $ws.Range("C3").Value = "南宁·第一届异次元动漫嘉年华"
    $ws.Range("F3").Value = 609
    $ws.Range("G3").Value = 5000

    # Row 4
    $ws.Range("C4").Value = "南宁·桂南动漫游戏嘉年华"
    $ws.Range("F4").Value = 1430
    $ws.Range("G4").Value = 6000

    # Row 5
    $ws.Range("C5").Value = "南宁·2024良牙动漫冬季盛典（冬典）"
    $ws.Range("F5").Value = 6948
    $ws.Range("G5").Value = 5500

    # Row 6
    $ws.Range("C6").Value = "南宁·第五届小蜜蜂动漫嘉年华"
    $ws.Range("F6").Value = 172
    $ws.Range("G6").Value = 5000

    # Row 7
    $ws.Range("C7").Value = "南宁·草莓动漫节"
    $ws.Range("F7").Value = 103
    $ws.Range("G7").Value = 6000
}
